# Pending task sheet updated:
# - add 5 new pending-task rows to the "tasks" sheet
# - update the sheet view (zoom + selection) to reflect where the user
#   was working after adding the rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tasks")
$ws.Activate()

# New pending tasks (rows 24/25 are left blank, matching the gap in the
# original edit, new data starts again at row 26)
$ws.Range("A26").Value = "addblog not working"
$ws.Range("A27").Value = "edit blog not working"
$ws.Range("A28").Value = "amenities add "
$ws.Range("A29").Value = "add_interior page not found need to be fixed"
$ws.Range("A30").Value = "edit_interior functionality not working"

# Reflect the resulting view state: zoomed in to 145% with A23 selected
$win = $excel.ActiveWindow
$win.Zoom = 145
$ws.Range("A23").Select()
